$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ticker values in columns B, C, E, F for rows 2-23
$ws.Range("B2").Value = "NSE:3PLAND"
$ws.Range("C2").Value = "NSE:ANANTRAJ"
$ws.Range("E2").Value = "NSE:INDIGO"
$ws.Range("F2").Value = "NSE:BRITANNIA"

$ws.Range("B3").Value = "NSE:ARVIND"
$ws.Range("C3").Value = "NSE:CCHHL"
$ws.Range("E3").ClearContents()
$ws.Range("F3").Value = "NSE:GAIL"

$ws.Range("B4").Value = "NSE:BIRLACORPN"
$ws.Range("C4").Value = "NSE:DPSCLTD"
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = "NSE:KOTAKBANK"

$ws.Range("B5").Value = "NSE:BRITANNIA"
$ws.Range("C5").Value = "NSE:FCSSOFT"
$ws.Range("F5").Value = "NSE:LTIM"

$ws.Range("B6").Value = "NSE:BSE"
$ws.Range("C6").Value = "NSE:GMBREW"
$ws.Range("F6").Value = "NSE:MOTHERSON"

$ws.Range("B7").Value = "NSE:CAMPUS"
$ws.Range("C7").Value = "NSE:GOCLCORP"
$ws.Range("F7").Value = "NSE:RELIANCE"

$ws.Range("B8").Value = "NSE:CEATLTD"
$ws.Range("C8").Value = "NSE:HAPPSTMNDS"

$ws.Range("B9").Value = "NSE:CIEINDIA"
$ws.Range("C9").Value = "NSE:HDFCMID150"

$ws.Range("B10").Value = "NSE:CMSINFO"
$ws.Range("C10").Value = "NSE:JAGRAN"

$ws.Range("B11").Value = "NSE:GAIL"
$ws.Range("C11").Value = "NSE:LIBERTSHOE"

$ws.Range("B12").Value = "NSE:GRASIM"
$ws.Range("C12").Value = "NSE:LIKHITHA"

$ws.Range("B13").Value = "NSE:GRPLTD"
$ws.Range("C13").Value = "NSE:MANUGRAPH"

$ws.Range("B14").Value = "NSE:INDUSTOWER"
$ws.Range("C14").Value = "NSE:MEGASTAR"

$ws.Range("B15").Value = "NSE:JKLAKSHMI"
$ws.Range("C15").Value = "NSE:MSUMI"

$ws.Range("B16").Value = "NSE:KESORAMIND"
$ws.Range("C16").Value = "NSE:NAHARCAP"

$ws.Range("B17").Value = "NSE:KOTAKBANK"
$ws.Range("C17").Value = "NSE:PLAZACABLE"

$ws.Range("B18").Value = "NSE:LTIM"
$ws.Range("C18").Value = "NSE:RAILTEL"

$ws.Range("B19").Value = "NSE:LTTS"
$ws.Range("C19").Value = "NSE:REMSONSIND"

$ws.Range("B20").Value = "NSE:MAZDOCK"
$ws.Range("C20").ClearContents()

$ws.Range("B21").Value = "NSE:NH"
$ws.Range("C21").ClearContents()

$ws.Range("B22").Value = "NSE:ORIENTCEM"
$ws.Range("C22").ClearContents()

$ws.Range("B23").Value = "NSE:RELIANCE"
$ws.Range("C23").ClearContents()

# Remove rows 24-42 entirely (no longer part of the data range)
$ws.Range("A24:A42").EntireRow.Delete()
